$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of sapling.id values for A2:A101 after adding site info to manu
# (values 149 and 1333 were relocated further down the list)
$values = @(
    27,32,36,42,55,63,88,92,103,118,125,150,163,171,176,207,228,235,271,280,
    289,308,328,340,381,412,424,426,427,429,430,433,434,437,447,450,455,460,492,498,
    509,512,514,525,544,606,622,728,731,740,762,791,814,822,891,905,927,1074,1123,1141,
    1145,1151,1154,1209,1249,1267,1268,1292,1307,1309,1310,1311,1318,1321,1323,1325,1326,1334,1336,1337,
    1338,1340,1343,1344,1385,1398,1404,1411,1437,1482,1492,113,565,9,97,142,
    149,1314,1331,1333
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
